# 17.1.2.xlsx — add the 2020 data point (year header + value) to the
# "Tax revenues" table on the single worksheet, and move the on-screen
# selection/scroll position the way the author last left the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently ends at column P (year 2019). Extend it one column
# to the right (Q) with year 2020 / value 70.3, inheriting the same
# look (borders, font, number format) as the preceding O:P columns by
# copying their formatting before writing the new numbers in.
$ws.Range("P4:P5").Copy($ws.Range("Q4:Q5"))

$ws.Cells.Item(4, 17).Value = 2020   # Q4 - year header
$ws.Cells.Item(5, 17).Value = 70.3   # Q5 - tax revenue % for 2020

# Restore the view: the saved file shows the grid scrolled right so
# column C is at the left edge, with Q8 as the active/selected cell.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("Q8").Select()
